$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# --- ALC row hunk 0 ---
$ws_ALC.Range("H17").Value = 1719.1111
$ws_ALC.Range("J17").Value = 1719.1111
$ws_ALC.Range("L17").Value = 5157.3333
$ws_ALC.Range("N17").Value = -5493.3333

# --- ALC row hunk 1 ---
$ws_ALC.Range("H34").Value = 11288.9
$ws_ALC.Range("J34").Value = 30003
$ws_ALC.Range("L34").Value = 30003
$ws_ALC.Range("N34").Value = -30409

# --- ALC row hunk 2 ---
$ws_ALC.Range("H36").Value = 11288.9
$ws_ALC.Range("J36").Value = 30003
$ws_ALC.Range("L36").Value = 30003
$ws_ALC.Range("N36").Value = -31433

# --- ALC row hunk 3 ---
$ws_ALC.Range("H86").Value = 4519.2
$ws_ALC.Range("I86").Value = 1299.1666
$ws_ALC.Range("J86").Value = 9349.25
$ws_ALC.Range("K86").Value = 1299.1666
$ws_ALC.Range("L86").Value = 9349.25
$ws_ALC.Range("M86").Value = -176.1666
$ws_ALC.Range("N86").Value = -11595.25

# --- ALC row hunk 4 ---
$ws_ALC.Range("H89").Value = 4519.2
$ws_ALC.Range("I89").Value = 1299.1666
$ws_ALC.Range("J89").Value = 9349.25
$ws_ALC.Range("K89").Value = 6495.833000000001
$ws_ALC.Range("L89").Value = 46746.25
$ws_ALC.Range("M89").Value = -879.8330000000005
$ws_ALC.Range("N89").Value = -57978.25

# --- ALC row hunk 5 ---
$ws_ALC.Range("H92").Value = 277.05884
$ws_ALC.Range("I92").Value = 377.5
$ws_ALC.Range("J92").Value = 133.57143
$ws_ALC.Range("K92").Value = 377.5
$ws_ALC.Range("L92").Value = 133.57143
$ws_ALC.Range("M92").Value = 870.5
$ws_ALC.Range("N92").Value = -2629.57143

# --- ALC row hunk 6 ---
$ws_ALC.Range("H98").Value = 2519.3333
$ws_ALC.Range("I98").Value = 2519.3333
$ws_ALC.Range("K98").Value = 2519.3333
$ws_ALC.Range("M98").Value = -1021.3333

# --- ALC row hunk 7 ---
$ws_ALC.Range("H99").Value = 2853.8
$ws_ALC.Range("I99").Value = 1900
$ws_ALC.Range("J99").Value = 3489.6667
$ws_ALC.Range("K99").Value = 5700
$ws_ALC.Range("L99").Value = 10469.0001
$ws_ALC.Range("M99").Value = -4202
$ws_ALC.Range("N99").Value = -13465.0001

# --- ALC row hunk 8 ---
$ws_ALC.Range("H101").Value = 1694.5
$ws_ALC.Range("J101").Value = 3995
$ws_ALC.Range("L101").Value = 11985
$ws_ALC.Range("N101").Value = -15229

# --- ALC row hunk 9 ---
$ws_ALC.Range("H103").Value = 0
$ws_ALC.Range("I103").Value = 0
$ws_ALC.Range("J103").Value = 0
$ws_ALC.Range("K103").Value = 0
$ws_ALC.Range("L103").Value = 0
$ws_ALC.Range("M103").ClearContents()
$ws_ALC.Range("N103").ClearContents()

# --- ALC row hunk 10 ---
$ws_ALC.Range("H122").Value = 2519.3333
$ws_ALC.Range("I122").Value = 2519.3333
$ws_ALC.Range("K122").Value = 7557.999899999999
$ws_ALC.Range("M122").Value = -5107.999899999999

# --- ALC row hunk 11 ---
$ws_ALC.Range("H132").Value = 0
$ws_ALC.Range("I132").Value = 0
$ws_ALC.Range("J132").Value = 0
$ws_ALC.Range("K132").Value = 0
$ws_ALC.Range("L132").Value = 0
$ws_ALC.Range("M132").ClearContents()
$ws_ALC.Range("N132").ClearContents()

# --- ALC row hunk 12 ---
$ws_ALC.Range("H137").Value = 1348.6316
$ws_ALC.Range("I137").Value = 1283.9375
$ws_ALC.Range("J137").Value = 1693.6666
$ws_ALC.Range("K137").Value = 3851.8125
$ws_ALC.Range("L137").Value = 5080.9998
$ws_ALC.Range("M137").Value = -1301.8125
$ws_ALC.Range("N137").Value = -10180.9998

# --- ARM row hunk 13 ---
$ws_ARM.Range("H2").Value = 623.7857
$ws_ARM.Range("J2").Value = 983.5
$ws_ARM.Range("L2").Value = 983.5
$ws_ARM.Range("N2").Value = -1209.5

# --- ARM row hunk 14 ---
$ws_ARM.Range("H15").Value = 0
$ws_ARM.Range("J15").Value = 0
$ws_ARM.Range("L15").Value = 0
$ws_ARM.Range("N15").ClearContents()

# --- ARM row hunk 15 ---
$ws_ARM.Range("H32").Value = 2924.0938
$ws_ARM.Range("I32").Value = 3056
$ws_ARM.Range("K32").Value = 3056
$ws_ARM.Range("M32").Value = -2769

# --- ARM row hunk 16 ---
$ws_ARM.Range("H61").Value = 3843.7778
$ws_ARM.Range("I61").Value = 3483
$ws_ARM.Range("J61").Value = 4294.75
$ws_ARM.Range("K61").Value = 3483
$ws_ARM.Range("L61").Value = 4294.75
$ws_ARM.Range("M61").Value = -3271
$ws_ARM.Range("N61").Value = -4718.75

# --- ARM row hunk 17 ---
$ws_ARM.Range("H74").Value = 15433234
$ws_ARM.Range("I74").Value = 9260962
$ws_ARM.Range("K74").Value = 9260962
$ws_ARM.Range("M74").Value = -9260088

# --- ARM row hunk 18 ---
$ws_ARM.Range("H77").Value = 15433234
$ws_ARM.Range("I77").Value = 9260962
$ws_ARM.Range("K77").Value = 46304810
$ws_ARM.Range("M77").Value = -46300442

# --- ARM row hunk 19 ---
$ws_ARM.Range("H88").Value = 2249.6
$ws_ARM.Range("I88").Value = 3000
$ws_ARM.Range("J88").Value = 1749.3334
$ws_ARM.Range("K88").Value = 3000
$ws_ARM.Range("L88").Value = 1749.3334
$ws_ARM.Range("M88").Value = -2594
$ws_ARM.Range("N88").Value = -2561.3334

# --- ARM row hunk 20 ---
$ws_ARM.Range("H91").Value = 2249.6
$ws_ARM.Range("I91").Value = 3000
$ws_ARM.Range("J91").Value = 1749.3334
$ws_ARM.Range("K91").Value = 3000
$ws_ARM.Range("L91").Value = 1749.3334
$ws_ARM.Range("M91").Value = -1596
$ws_ARM.Range("N91").Value = -4557.3334

# --- ARM row hunk 21 ---
$ws_ARM.Range("H110").Value = 1027
$ws_ARM.Range("I110").Value = 969.3333
$ws_ARM.Range("J110").Value = 1200
$ws_ARM.Range("K110").Value = 969.3333
$ws_ARM.Range("L110").Value = 1200
$ws_ARM.Range("M110").Value = 1075.6667
$ws_ARM.Range("N110").Value = -5290

# --- ARM row hunk 22 ---
$ws_ARM.Range("H116").Value = 623.7857
$ws_ARM.Range("J116").Value = 983.5
$ws_ARM.Range("L116").Value = 983.5
$ws_ARM.Range("N116").Value = -5571.5

# --- ARM row hunk 23 ---
$ws_ARM.Range("H136").Value = 3843.7778
$ws_ARM.Range("I136").Value = 3483
$ws_ARM.Range("J136").Value = 4294.75
$ws_ARM.Range("K136").Value = 10449
$ws_ARM.Range("L136").Value = 12884.25
$ws_ARM.Range("M136").Value = -7899
$ws_ARM.Range("N136").Value = -17984.25

# --- BSM row hunk 24 ---
$ws_BSM.Range("H3").Value = 623.7857
$ws_BSM.Range("J3").Value = 983.5
$ws_BSM.Range("L3").Value = 983.5
$ws_BSM.Range("N3").Value = -1211.5

# --- BSM row hunk 25 ---
$ws_BSM.Range("H86").Value = 5975
$ws_BSM.Range("I86").Value = 7938
$ws_BSM.Range("K86").Value = 7938
$ws_BSM.Range("M86").Value = -6815

# --- BSM row hunk 26 ---
$ws_BSM.Range("H89").Value = 5975
$ws_BSM.Range("I89").Value = 7938
$ws_BSM.Range("K89").Value = 39690
$ws_BSM.Range("M89").Value = -34074

# --- BSM row hunk 27 ---
$ws_BSM.Range("H99").Value = 4098.647
$ws_BSM.Range("I99").Value = 4079.8125
$ws_BSM.Range("K99").Value = 4079.8125
$ws_BSM.Range("M99").Value = -2581.8125

# --- CRP row hunk 28 ---
$ws_CRP.Range("H16").Value = 2059.9167
$ws_CRP.Range("I16").Value = 1613.375
$ws_CRP.Range("J16").Value = 2953
$ws_CRP.Range("K16").Value = 1613.375
$ws_CRP.Range("L16").Value = 2953
$ws_CRP.Range("M16").Value = -1326.375
$ws_CRP.Range("N16").Value = -3527

# --- CRP row hunk 29 ---
$ws_CRP.Range("H18").Value = 34880.5
$ws_CRP.Range("J18").Value = 34880.5
$ws_CRP.Range("L18").Value = 34880.5
$ws_CRP.Range("N18").Value = -35340.5

# --- CRP row hunk 30 ---
$ws_CRP.Range("H102").Value = 19499.334
$ws_CRP.Range("J102").Value = 19499.334
$ws_CRP.Range("L102").Value = 19499.334
$ws_CRP.Range("N102").Value = -24367.334

# --- CRP row hunk 31 ---
$ws_CRP.Range("H109").Value = 91768.78999999999
$ws_CRP.Range("J109").Value = 91768.78999999999
$ws_CRP.Range("L109").Value = 91768.78999999999
$ws_CRP.Range("N109").Value = -93848.78999999999

# --- CRP row hunk 32 ---
$ws_CRP.Range("H113").Value = 2059.9167
$ws_CRP.Range("I113").Value = 1613.375
$ws_CRP.Range("J113").Value = 2953
$ws_CRP.Range("K113").Value = 1613.375
$ws_CRP.Range("L113").Value = 2953
$ws_CRP.Range("M113").Value = 556.625
$ws_CRP.Range("N113").Value = -7293

# --- CRP row hunk 33 ---
$ws_CRP.Range("H122").Value = 2499
$ws_CRP.Range("I122").Value = 2499
$ws_CRP.Range("J122").Value = 2499
$ws_CRP.Range("K122").Value = 7497
$ws_CRP.Range("L122").Value = 7497
$ws_CRP.Range("N122").Value = -12397
$ws_CRP.Range("M122").Value = -5047

# --- CUL row hunk 34 ---
$ws_CUL.Range("H26").Value = 0
$ws_CUL.Range("I26").Value = 0
$ws_CUL.Range("J26").Value = 0
$ws_CUL.Range("K26").Value = 0
$ws_CUL.Range("L26").Value = 0
$ws_CUL.Range("M26").ClearContents()
$ws_CUL.Range("N26").ClearContents()

# --- CUL row hunk 35 ---
$ws_CUL.Range("H122").Value = 1619.5385
$ws_CUL.Range("J122").Value = 2458.1667
$ws_CUL.Range("L122").Value = 22123.5003
$ws_CUL.Range("N122").Value = -27023.5003

# --- GSM row hunk 36 ---
$ws_GSM.Range("H107").Value = 1040.4375
$ws_GSM.Range("I107").Value = 834.3
$ws_GSM.Range("K107").Value = 834.3
$ws_GSM.Range("M107").Value = 1085.7

# --- GSM row hunk 37 ---
$ws_GSM.Range("H132").Value = 1100.9166
$ws_GSM.Range("I132").Value = 1094.2222
$ws_GSM.Range("J132").Value = 1121
$ws_GSM.Range("K132").Value = 3282.6666
$ws_GSM.Range("L132").Value = 3363
$ws_GSM.Range("M132").Value = -752.6665999999996
$ws_GSM.Range("N132").Value = -8423

# --- LTW row hunk 38 ---
$ws_LTW.Range("H109").Value = 43998
$ws_LTW.Range("J109").Value = 43998
$ws_LTW.Range("L109").Value = 43998
$ws_LTW.Range("N109").Value = -46772

# --- WVR row hunk 39 ---
$ws_WVR.Range("H2").Value = 26570.857
$ws_WVR.Range("I2").Value = 9999.5
$ws_WVR.Range("J2").Value = 33199.4
$ws_WVR.Range("K2").Value = 9999.5
$ws_WVR.Range("L2").Value = 33199.4
$ws_WVR.Range("M2").Value = -9887.5
$ws_WVR.Range("N2").Value = -33423.4

# --- WVR row hunk 40 ---
$ws_WVR.Range("H116").Value = 66589
$ws_WVR.Range("J116").Value = 66589
$ws_WVR.Range("L116").Value = 66589
$ws_WVR.Range("N116").Value = -75767

# --- WVR row hunk 41 ---
$ws_WVR.Range("H122").Value = 4924.4287
$ws_WVR.Range("I122").Value = 4924.4287
$ws_WVR.Range("K122").Value = 14773.2861
$ws_WVR.Range("M122").Value = -12323.2861
